$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.310.12"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "1.588.79"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.38%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "210.03"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.506"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.74%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "19.45"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.71%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0844"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "1.812.07"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.07"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.584.81"
$ws.Range("E14").Value = "  -0.67%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.519"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "26.315.01"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "0.0₃0728"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.49"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +6.02%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "211.23"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E23").Value = "  -3.12%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "8.94"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "144.65"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E26").Value = "  -0.32%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.05"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("E28").Value = "  -0.68%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "15.24"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0505"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "1.306.14"
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.611"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.02%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.44"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  -9.53%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.805"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("E41").Value = "  -0.33%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.60"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("E43").Value = "  -0.51%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.12"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").Value = "1.724.43"
$ws.Range("E46").Value = "  -0.46%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "87.84"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("E48").Value = "  -5.25%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0506"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.31%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0980"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.85%  "
$ws.Range("E51").Value = "  -0.31%  "
